# Nexial "redis-showcase" workbook update:
# Adds a new "localdb" command-type category to the hidden '#system' sheet.
# This requires:
#   1. Inserting a new column at N (pushing existing N:AC categories to O:AD)
#      and filling the freed column N with the "localdb" category header and
#      its six function names.
#   2. Inserting "localdb" into the alphabetically-sorted "target" category
#      list in column A (between "json" and "macro"), which pushes A14:A29
#      down to A15:A30.
#   3. Updating / adding the matching named ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1. Shift columns N:AC (14:29) right by one, into O:AD (15:30), using a
#    bulk array copy so unrelated columns/rows are left untouched.
# ---------------------------------------------------------------------
$srcRange = $ws.Range("N1:AC127")
$shiftedValues = $srcRange.Value()
$dstRange = $ws.Range("O1:AD127")
$dstRange.Value = $shiftedValues

# clear out the now-stale copy left behind in column N
$ws.Range("N1:N127").ClearContents()

# populate the freed column N with the new "localdb" category
$ws.Range("N1").Value = "localdb"
$ws.Range("N2").Value = "cloneTable(var,source,target)"
$ws.Range("N3").Value = "dropTables(var,tables)"
$ws.Range("N4").Value = "exportCSV(sql,output)"
$ws.Range("N5").Value = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value = "purge(var)"
$ws.Range("N7").Value = "runSQLs(var,sqls)"

# ---------------------------------------------------------------------
# 2. Insert "localdb" into column A's "target" list (row 14), pushing the
#    existing entries (macro, mail, number, ... xml) down by one row.
# ---------------------------------------------------------------------
$srcA = $ws.Range("A14:A29")
$shiftedA = $srcA.Value()
$dstA = $ws.Range("A15:A30")
$dstA.Value = $shiftedA

$ws.Range("A14").Value = "localdb"

# ---------------------------------------------------------------------
# 3. Update the named ranges that refer to the shifted columns, and add
#    the new "localdb" named range.
# ---------------------------------------------------------------------
$names = $wb.Names

$names.Item("mail").RefersTo        = "='#system'!`$P`$2:`$P`$2"
$names.Item("number").RefersTo      = "='#system'!`$Q`$2:`$Q`$16"
$names.Item("pdf").RefersTo         = "='#system'!`$R`$2:`$R`$16"
$names.Item("rdbms").RefersTo       = "='#system'!`$S`$2:`$S`$7"
$names.Item("redis").RefersTo       = "='#system'!`$T`$2:`$T`$10"
$names.Item("sms").RefersTo         = "='#system'!`$U`$2:`$U`$2"
$names.Item("sound").RefersTo       = "='#system'!`$V`$2:`$V`$5"
$names.Item("ssh").RefersTo         = "='#system'!`$W`$2:`$W`$9"
$names.Item("step").RefersTo        = "='#system'!`$X`$2:`$X`$4"
$names.Item("target").RefersTo      = "='#system'!`$A`$2:`$A`$30"
$names.Item("web").RefersTo         = "='#system'!`$Y`$2:`$Y`$127"
$names.Item("webalert").RefersTo    = "='#system'!`$Z`$2:`$Z`$8"
$names.Item("webcookie").RefersTo   = "='#system'!`$AA`$2:`$AA`$8"
$names.Item("ws").RefersTo          = "='#system'!`$AB`$2:`$AB`$17"
$names.Item("ws.async").RefersTo    = "='#system'!`$AC`$2:`$AC`$8"
$names.Item("xml").RefersTo         = "='#system'!`$AD`$2:`$AD`$21"
$names.Item("macro").RefersTo       = "='#system'!`$O`$2:`$O`$4"

$names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
